$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row  # xlCellTypeLastCell -> 11
$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
